$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, shifting rows 84:150 down to 85:151
$ws.Rows.Item(84).Insert()

# Populate the new row 84 with fresh data
$ws.Cells.Item(84, 1).Value = 10
$ws.Cells.Item(84, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(84, 3).Value = "La Araucanía"
$ws.Cells.Item(84, 4).Value = 44658
$ws.Cells.Item(84, 4).NumberFormat = $ws.Cells.Item(85, 4).NumberFormat
$ws.Cells.Item(84, 5).Value = 9
$ws.Cells.Item(84, 6).Value = 100114007
$ws.Cells.Item(84, 7).Value = "Jengibre"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 65
$ws.Cells.Item(84, 11).Value = 25000
$ws.Cells.Item(84, 12).Value = 25000
$ws.Cells.Item(84, 13).Value = 25000
$ws.Cells.Item(84, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(84, 15).Value = "Perú"
$ws.Cells.Item(84, 16).Value = 1923
$ws.Cells.Item(84, 17).Value = 13
$ws.Cells.Item(84, 18).Value = "Hortaliza"
